$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04279549268760489
$ws.Range("C3").Value = 0.06772261264820557
$ws.Range("D3").Value = 0.02492711996060069
$ws.Range("E3").Value = -0.01047939650004155
$ws.Range("F3").Value = 0.03540651646064223
$ws.Range("G3").Value = 0.2283793832111516
$ws.Range("H3").Value = 0.7716206167888483
$ws.Range("C4").Value = 0.07653187936772365
$ws.Range("D4").Value = 0.03373638668011877
$ws.Range("E4").Value = -0.007267482430366687
$ws.Range("F4").Value = 0.04100386911048545
$ws.Range("G4").Value = 0.1505547741752415
$ws.Range("H4").Value = 0.8494452258247586
$ws.Range("C5").Value = 0.08049062233732643
$ws.Range("D5").Value = 0.03769512964972154
$ws.Range("E5").Value = -0.007701953802140163
$ws.Range("F5").Value = 0.04539708345186169
$ws.Range("G5").Value = 0.14504884081603
$ws.Range("H5").Value = 0.8549511591839699
$ws.Range("D6").Value = 0.04248059810190562
$ws.Range("E6").Value = -0.01173544055097419
$ws.Range("F6").Value = 0.05421603865287981
$ws.Range("G6").Value = 0.17794052070766
$ws.Range("H6").Value = 0.8220594792923401
$ws.Range("C7").Value = 0.09151276883940666
$ws.Range("D7").Value = 0.04871727615180176
$ws.Range("E7").Value = -0.007802166799712454
$ws.Range("F7").Value = 0.05651944295151422
$ws.Range("G7").Value = 0.1212993087375843
$ws.Range("H7").Value = 0.8787006912624158
$ws.Range("C8").Value = 0.07887698887094986
$ws.Range("D8").Value = 0.03608149618334498
$ws.Range("E8").Value = -0.008945727429587139
$ws.Range("F8").Value = 0.04502722361293211
$ws.Range("G8").Value = 0.1657446416546652
$ws.Range("H8").Value = 0.8342553583453349
$ws.Range("D9").Value = 0.03852701461435768
$ws.Range("E9").Value = -0.009197806527930343
$ws.Range("F9").Value = 0.04772482114228803
$ws.Range("G9").Value = 0.1615843629218577
$ws.Range("H9").Value = 0.8384156370781424
$ws.Range("C10").Value = 0.08841811769411888
$ws.Range("D10").Value = 0.04562262500651399
$ws.Range("E10").Value = -0.008760123667387685
$ws.Range("F10").Value = 0.05438274867390167
$ws.Range("G10").Value = 0.1387349568141107
$ws.Range("H10").Value = 0.8612650431858894
$ws.Range("C11").Value = 0.08655089758402457
$ws.Range("D11").Value = 0.04375540489641968
$ws.Range("E11").Value = -0.01134134369417696
$ws.Range("F11").Value = 0.05509674859059664
$ws.Range("G11").Value = 0.1707054387649266
$ws.Range("H11").Value = 0.8292945612350734
$ws.Range("C12").Value = 0.0919005225823699
$ws.Range("D12").Value = 0.04910502989476501
$ws.Range("E12").Value = -0.007595670867927059
$ws.Range("F12").Value = 0.05670070076269207
$ws.Range("G12").Value = 0.1181352955896793
$ws.Range("H12").Value = 0.8818647044103207
$ws.Range("C13").Value = 0.09303340607688435
$ws.Range("D13").Value = 0.05023791338927946
$ws.Range("E13").Value = -0.00777813616395686
$ws.Range("F13").Value = 0.05801604955323631
$ws.Range("G13").Value = 0.1182192024898683
$ws.Range("H13").Value = 0.8817807975101317
$ws.Range("C14").Value = 0.09761224933057054
$ws.Range("D14").Value = 0.05481675664296565
$ws.Range("E14").Value = -0.006745914997421222
$ws.Range("F14").Value = 0.06156267164038688
$ws.Range("G14").Value = 0.09875647161593339
$ws.Range("H14").Value = 0.9012435283840665
$ws.Range("C15").Value = 0.09217957736456678
$ws.Range("D15").Value = 0.04938408467696189
$ws.Range("E15").Value = -0.01130594266591633
$ws.Range("F15").Value = 0.06069002734287823
$ws.Range("G15").Value = 0.1570357710929552
$ws.Range("H15").Value = 0.8429642289070448
$ws.Range("C16").Value = 0.09630552684121836
$ws.Range("D16").Value = 0.05351003415361347
$ws.Range("E16").Value = -0.008642120971093972
$ws.Range("F16").Value = 0.06215215512470745
$ws.Range("G16").Value = 0.1220737247090307
$ws.Range("H16").Value = 0.8779262752909695
$ws.Range("C17").Value = 0.09852677163207846
$ws.Range("D17").Value = 0.05573127894447357
$ws.Range("E17").Value = -0.008912367949870338
$ws.Range("F17").Value = 0.0646436468943439
$ws.Range("G17").Value = 0.1211643666224444
$ws.Range("H17").Value = 0.8788356333775556
$ws.Range("D18").Value = 0.06127634886732191
$ws.Range("E18").Value = -0.008534332691252207
$ws.Range("F18").Value = 0.06981068155857413
$ws.Range("G18").Value = 0.1089326841404094
$ws.Range("H18").Value = 0.8910673158595905
$ws.Range("C19").Value = 0.1068333978350594
$ws.Range("D19").Value = 0.06403790514745453
$ws.Range("E19").Value = -0.009089977724716526
$ws.Range("F19").Value = 0.07312788287217106
$ws.Range("G19").Value = 0.1105596479733825
$ws.Range("H19").Value = 0.8894403520266175
$ws.Range("C20").Value = 0.1063066138915367
$ws.Range("D20").Value = 0.06351112120393182
$ws.Range("E20").Value = -0.01193157164180236
$ws.Range("F20").Value = 0.0754426928457342
$ws.Range("G20").Value = 0.1365570481397796
$ws.Range("H20").Value = 0.8634429518602205
$ws.Range("C21").Value = 0.1094216564513594
$ws.Range("D21").Value = 0.06662616376375452
$ws.Range("E21").Value = -0.01007276215778308
$ws.Range("F21").Value = 0.07669892592153758
$ws.Range("G21").Value = 0.1160835104253735
$ws.Range("H21").Value = 0.8839164895746265
$ws.Range("C22").Value = 0.1160763913934513
$ws.Range("D22").Value = 0.07328089870584645
$ws.Range("E22").Value = -0.008594126959789552
$ws.Range("F22").Value = 0.081875025665636
$ws.Range("G22").Value = 0.09499510839205373
$ws.Range("H22").Value = 0.9050048916079463
$ws.Range("D23").Value = 0.07551611476172579
$ws.Range("E23").Value = -0.01006122870785813
$ws.Range("F23").Value = 0.08557734346958391
$ws.Range("G23").Value = 0.1052005323667017
$ws.Range("H23").Value = 0.8947994676332983
$ws.Range("C24").Value = 0.1226522598312383
$ws.Range("D24").Value = 0.07985676714363343
$ws.Range("E24").Value = -0.009450035996411895
$ws.Range("F24").Value = 0.08930680314004531
$ws.Range("G24").Value = 0.09568993984664002
$ws.Range("H24").Value = 0.90431006015336
$ws.Range("C25").Value = 0.1184248271280611
$ws.Range("D25").Value = 0.07562933444045622
$ws.Range("E25").Value = -0.01198187879081628
$ws.Range("F25").Value = 0.08761121323127249
$ws.Range("G25").Value = 0.1203083321096087
$ws.Range("H25").Value = 0.8796916678903913
$ws.Range("C26").Value = 0.1181602872777506
$ws.Range("D26").Value = 0.07536479459014572
$ws.Range("E26").Value = -0.01195920358542672
$ws.Range("F26").Value = 0.08732399817557245
$ws.Range("G26").Value = 0.1204554584592838
$ws.Range("H26").Value = 0.8795445415407163
$ws.Range("C27").Value = 0.1188512226001171
$ws.Range("D27").Value = 0.07605572991251217
$ws.Range("E27").Value = -0.01201565855867775
$ws.Range("F27").Value = 0.08807138847118992
$ws.Range("G27").Value = 0.1200520838135236
$ws.Range("H27").Value = 0.8799479161864764
$ws.Range("C28").Value = 0.1184472032807328
$ws.Range("D28").Value = 0.07565171059312789
$ws.Range("E28").Value = -0.01289169010586391
$ws.Range("F28").Value = 0.0885434006989918
$ws.Range("G28").Value = 0.1270930010864326
$ws.Range("H28").Value = 0.8729069989135675
$ws.Range("C29").Value = 0.1193307475487244
$ws.Range("D29").Value = 0.07653525486111951
$ws.Range("E29").Value = -0.01201674115314846
$ws.Range("F29").Value = 0.08855199601426796
$ws.Range("G29").Value = 0.1194878397761348
$ws.Range("H29").Value = 0.8805121602238653
$ws.Range("C30").Value = 0.1182901152261737
$ws.Range("D30").Value = 0.07549462253856876
$ws.Range("E30").Value = -0.0125874609386457
$ws.Range("F30").Value = 0.08808208347721445
$ws.Range("G30").Value = 0.1250374282677551
$ws.Range("H30").Value = 0.8749625717322449
$ws.Range("C31").Value = 0.1217573532574061
$ws.Range("D31").Value = 0.07896186056980117
$ws.Range("E31").Value = -0.01213877532354349
$ws.Range("F31").Value = 0.09110063589334465
$ws.Range("G31").Value = 0.1175788895002706
$ws.Range("H31").Value = 0.8824211104997294
$ws.Range("C32").Value = 0.1258695347424495
$ws.Range("D32").Value = 0.08307404205484462
$ws.Range("E32").Value = -0.00999570396931908
$ws.Range("F32").Value = 0.09306974602416368
$ws.Range("G32").Value = 0.0969840423726006
$ws.Range("H32").Value = 0.9030159576273994
$ws.Range("C33").Value = 0.122356092604443
$ws.Range("D33").Value = 0.07956059991683812
$ws.Range("E33").Value = -0.01158899952857768
$ws.Range("F33").Value = 0.09114959944541581
$ws.Range("G33").Value = 0.1128008328350987
$ws.Range("H33").Value = 0.8871991671649013
$ws.Range("C34").Value = 0.1243886590084101
$ws.Range("D34").Value = 0.08159316632080524
$ws.Range("E34").Value = -0.010450814751991
$ws.Range("F34").Value = 0.09204398107279625
$ws.Range("G34").Value = 0.1019643452908229
$ws.Range("H34").Value = 0.8980356547091771
$ws.Range("C35").Value = 0.1268083105054177
$ws.Range("D35").Value = 0.08401281781781284
$ws.Range("E35").Value = -0.01043916346306427
$ws.Range("F35").Value = 0.09445198128087709
$ws.Range("G35").Value = 0.09952378237980138
$ws.Range("H35").Value = 0.9004762176201986
$ws.Range("C36").Value = 0.1272850117438769
$ws.Range("D36").Value = 0.08448951905627203
$ws.Range("E36").Value = -0.01009605021245804
$ws.Range("F36").Value = 0.09458556926873005
$ws.Range("G36").Value = 0.09644530016343852
$ws.Range("H36").Value = 0.9035546998365615
$ws.Range("C37").Value = 0.1281204503424168
$ws.Range("D37").Value = 0.08532495765481188
$ws.Range("E37").Value = -0.009688529023479562
$ws.Range("F37").Value = 0.09501348667829145
$ws.Range("G37").Value = 0.09253431233908596
$ws.Range("H37").Value = 0.907465687660914
$ws.Range("D38").Value = 0.08460250157429472
$ws.Range("E38").Value = -0.008829469244627762
$ws.Range("F38").Value = 0.09343197081892248
$ws.Range("G38").Value = 0.08634211721584109
$ws.Range("H38").Value = 0.9136578827841589
$ws.Range("C39").Value = 0.1254792814685551
$ws.Range("D39").Value = 0.0826837887809502
$ws.Range("E39").Value = -0.01007294757108519
$ws.Range("F39").Value = 0.09275673635203538
$ws.Range("G39").Value = 0.09795758565801011
$ws.Range("H39").Value = 0.9020424143419898
$ws.Range("C40").Value = 0.125676431583777
$ws.Range("D40").Value = 0.08288093889617207
$ws.Range("E40").Value = -0.01054772225316979
$ws.Range("F40").Value = 0.09342866114934187
$ws.Range("G40").Value = 0.1014434423280296
$ws.Range("H40").Value = 0.8985565576719704
$ws.Range("C41").Value = 0.1267307432342762
$ws.Range("D41").Value = 0.08393525054667128
$ws.Range("E41").Value = -0.009344758275452674
$ws.Range("F41").Value = 0.09328000882212395
$ws.Range("G41").Value = 0.0910575345478503
$ws.Range("H41").Value = 0.9089424654521497
$ws.Range("C42").Value = 0.1272752828730058
$ws.Range("D42").Value = 0.08447979018540094
$ws.Range("E42").Value = -0.008995464874863166
$ws.Range("F42").Value = 0.09347525506026412
$ws.Range("G42").Value = 0.08778570971842556
$ws.Range("H42").Value = 0.9122142902815745
$ws.Range("C43").Value = 0.1300061819994451
$ws.Range("D43").Value = 0.08721068931184023
$ws.Range("E43").Value = -0.008081504117990253
$ws.Range("F43").Value = 0.09529219342983047
$ws.Range("G43").Value = 0.07817756653477298
$ws.Range("H43").Value = 0.9218224334652271
